$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells that look like plain numbers to remain text,
# matching the source data's inlineStr string representation (e.g. '22.00', '0.619').
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '51.433.55'
$ws.Range("E2").Value = '  -0.98%  '

$ws.Range("D3").Value = '2.780.48'
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '352.97'
$ws.Range("E5").Value = '  -1.95%  '

$ws.Range("D6").Value = '108.01'
$ws.Range("E6").Value = '  -1.44%  '

$ws.Range("E7").Value = '  -1.60%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '0.619'
$ws.Range("E9").Value = '  +4.78%  '

$ws.Range("D10").Value = '39.17'
$ws.Range("E10").Value = '  -2.40%  '

$ws.Range("E11").Value = '  +1.53%  '

$ws.Range("D12").Value = '0.0833'
$ws.Range("E12").Value = '  -1.74%  '

$ws.Range("D13").Value = '19.90'
$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("D14").Value = '7.75'
$ws.Range("E14").Value = '  +2.47%  '

$ws.Range("D15").Value = '3.215.61'
$ws.Range("E15").Value = '  -0.45%  '

$ws.Range("D16").Value = '2.790.49'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '0.924'
$ws.Range("E17").Value = '  -1.32%  '

$ws.Range("D18").Value = '51.429.33'
$ws.Range("E18").Value = '  -0.90%  '

$ws.Range("D19").Value = '7.71'
$ws.Range("E19").Value = '  +3.31%  '

$ws.Range("D20").Value = '3.11'
$ws.Range("E20").Value = '  +0.68%  '

$ws.Range("D21").Value = '13.36'
$ws.Range("E21").Value = '  +1.90%  '

$ws.Range("D22").Value = '0.0₃0966'
$ws.Range("E22").Value = '  -0.98%  '

$ws.Range("D23").Value = '70.50'
$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("D24").Value = '266.50'
$ws.Range("E24").Value = '  -1.22%  '

$ws.Range("D25").Value = '2.77'
$ws.Range("E25").Value = '  +0.79%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").Value = '25.84'
$ws.Range("E27").Value = '  -2.62%  '

$ws.Range("E28").Value = '  +1.37%  '

$ws.Range("D29").Value = '10.27'
$ws.Range("E29").Value = '  -0.30%  '

$ws.Range("D30").Value = '37.06'
$ws.Range("E30").Value = '  +7.94%  '

$ws.Range("E31").Value = '  -2.18%  '

$ws.Range("D32").Value = '6.22'
$ws.Range("E32").Value = '  +8.37%  '

$ws.Range("D33").Value = '51.80'
$ws.Range("E33").Value = '  -0.28%  '

$ws.Range("D34").Value = '5.68'
$ws.Range("E34").Value = '  +8.44%  '

$ws.Range("D35").Value = '0.0443'
$ws.Range("E35").Value = '  -5.63%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.0845'
$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").Value = '18.50'
$ws.Range("E38").Value = '  -2.81%  '

$ws.Range("E39").Value = '  -3.06%  '

$ws.Range("E40").Value = '  -1.80%  '

$ws.Range("E41").Value = '  -0.93%  '

$ws.Range("D42").Value = '2.49'
$ws.Range("E42").Value = '  -5.00%  '

$ws.Range("D43").Value = '120.06'
$ws.Range("E43").Value = '  +0.61%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '22.00'
$ws.Range("E44").Value = '  +0.81%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '2.18'
$ws.Range("E45").Value = '  -2.71%  '

$ws.Range("D46").Value = '2.130.48'
$ws.Range("E46").Value = '  +2.26%  '

$ws.Range("D47").Value = '3.35'
$ws.Range("E47").Value = '  +3.02%  '

$ws.Range("D48").Value = '2.33'
$ws.Range("E48").Value = '  +5.27%  '

$ws.Range("D49").Value = '0.228'
$ws.Range("E49").Value = '  +18.83%  '

$ws.Range("D50").Value = '5.48'
$ws.Range("E50").Value = '  -5.41%  '

$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").Value = '0.892'
$ws.Range("E51").Value = '  -6.92%  '

# Restore default (unstyled) appearance now that the text values are locked in.
$ws.Range("D2:D51").Style = "Normal"
